$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.742.87'
$ws.Range('E2').Value = '  -3.52%  '
$ws.Range('D3').Value = '3.123.14'
$ws.Range('E3').Value = '  -2.25%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '604.30'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.31%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '149.39'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.32%  '
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('D8').Value = '3.119.87'
$ws.Range('E8').Value = '  -2.38%  '
$ws.Range('E9').Value = '  -3.34%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.153'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.62%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.58'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.36%  '
$ws.Range('E12').Value = '  -5.03%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000259'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.78'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.81%  '
$ws.Range('D15').Value = '3.607.24'
$ws.Range('E15').Value = '  -3.02%  '
$ws.Range('D16').Value = '63.901.01'
$ws.Range('E16').Value = '  -3.50%  '
$ws.Range('E17').Value = '  +0.07%  '
$ws.Range('D18').Value = '3.116.10'
$ws.Range('E18').Value = '  -2.54%  '
$ws.Range('E19').Value = '  -4.62%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '483.07'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.85%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.61'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.46%  '
$ws.Range('E22').Value = '  -2.72%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.77'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.72%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.78'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -5.47%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.03'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.13%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.92'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.17%  '
$ws.Range('E28').Value = '  -4.95%  '
$ws.Range('B29').Value = 'Hedera'
$ws.Range('C29').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.127'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.62%  '
$ws.Range('B30').Value = 'ImmutableX'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.25'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.42%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.00'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.44%  '
$ws.Range('E32').Value = '  -0.09%  '
$ws.Range('E33').Value = '  -7.67%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '26.71'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.17%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.10'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.70%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.10'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.78%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '54.45'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.77%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.26'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +7.71%  '
$ws.Range('D39').Value = '0.0₃0747'
$ws.Range('E39').Value = '  -2.68%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '449.29'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -10.07%  '
$ws.Range('E41').Value = '  -4.62%  '
$ws.Range('E42').Value = '  -4.98%  '
$ws.Range('E43').Value = '  -2.58%  '
$ws.Range('D44').Value = '2.876.68'
$ws.Range('E44').Value = '  -1.14%  '
$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.272'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -8.14%  '
$ws.Range('B46').Value = 'Fetch.AI'
$ws.Range('C46').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.33'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.88%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '26.70'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.05%  '
$ws.Range('E48').Value = '  -0.05%  '
$ws.Range('E49').Value = '  -0.80%  '
$ws.Range('E50').Value = '  -3.66%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '119.54'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.34%  '
